# -----------------------------------------------------------------------
# Add 2022-Q3 data:
#   1. Shift the "总计" (summary) sheet rows down by one and insert the
#      new 2022-Q3 summary row at the top.
#   2. Insert a brand-new "2022-Q3" worksheet (positioned right after
#      "总计", before "2022-Q2") with the per-fund holdings detail.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet before "2022-Q2".
#    NOTE: worksheet object references become stale/rebound once a new
#    sheet is inserted, so we re-fetch every sheet handle by name right
#    after this call and before doing anything else.
# ---------------------------------------------------------------------
$insertBefore = $wb.Worksheets.Item("2022-Q2")
$added = $wb.Worksheets.Add($insertBefore)
$added.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Update the "总计" worksheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Shift existing rows 2-8 down to rows 3-9, preserving values & styles.
$summary.Range("A2:D8").Copy($summary.Range("A3:D9"))

# Write the new top row (2022-Q3 summary figures).
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 13
$summary.Range("D2").Value = 1.4

# ---------------------------------------------------------------------
# 3) Populate the new "2022-Q3" worksheet
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Item("2022-Q3")
$refSheet = $wb.Worksheets.Item("2022-Q2")

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holdings detail rows
$rows = @(
    @(0,  "483003", "工银精选平衡混合",               "15.78", "65.64", "2.41", "0.3803", 9),
    @(1,  "000893", "工银创新动力股票",                "11.92", "81.96", "3.19", "0.3802", 7),
    @(2,  "013233", "华夏中证500指数智选增强A",        "21.06", "93.71", "1.01", "0.2127", 10),
    @(3,  "007994", "华夏中证500指数增强A",            "19.56", "93.37", "1.00", "0.1956", 10),
    @(4,  "515450", "南方标普中国A股大盘红利低波50ETF", "2.17",  "99.66", "2.53", "0.0549", 10),
    @(5,  "007995", "华夏中证500指数增强C",            "5.27",  "93.37", "1.00", "0.0527", 10),
    @(6,  "013234", "华夏中证500指数智选增强C",        "3.92",  "93.71", "1.01", "0.0396", 10),
    @(7,  "011376", "华宝安享混合",                   "6.06",  "20.12", "0.54", "0.0327", 6),
    @(8,  "014133", "工银中证500六个月持有指数增强A",   "1.64",  "93.75", "1.87", "0.0307", 5),
    @(9,  "014134", "工银中证500六个月持有指数增强C",   "0.90",  "93.75", "1.87", "0.0168", 5),
    @(10, "005053", "银河量化价值混合A",               "0.10",  "78.55", "1.67", "0.0017", 8),
    @(11, "005126", "银河量化稳进混合",                "0.13",  "55.69", "1.10", "0.0014", 6),
    @(12, "013026", "银河量化价值混合C",               "0.00",  "78.55", "1.67", 0,        8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    if ($r -eq 14) {
        $newSheet.Cells.Item($r, 7).Value = $row[6]
    } else {
        $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    }
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 4) Apply formatting (bold/centered/bordered style) to header row and
#    the index column, matching the style used in the other sheets.
# ---------------------------------------------------------------------
$refSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$refSheet.Range("A2").Copy()
$newSheet.Range("A2:A14").PasteSpecial(-4122)

$newSheet.Range("A1").Select()
